$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update goal pose / orientation values for experimentA
$ws.Range("B12").Value = 154
$ws.Range("B13").Value = -339
$ws.Range("B15").Value = 0.34
$ws.Range("B16").Value = 0.62
$ws.Range("B17").Value = 0.62
$ws.Range("B18").Value = 0.34

# Move active selection to B18
$ws.Range("B18").Select()
